# Borehole survey: convert INCLINATION column (D) from "degrees from
# horizontal" (negative, e.g. -90) to "degrees from vertical" (e.g. 0)
# by adding 90 to every value in D2:D33.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 33; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $cell.Value2 = $cell.Value2 + 90
}

# Move the active selection to K14, matching the saved workbook state.
$ws.Range("K14").Select()
